$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.138.65'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '3.796.98'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  +0.09%  '
$c = $ws.Range('D5')
$c.Value = "'601.24"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$c = $ws.Range('D6')
$c.Value = "'163.74"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.27%  '
$ws.Range('D7').Value = '3.794.88'
$ws.Range('E7').Value = '  +1.18%  '
$ws.Range('E8').Value = '  -0.05%  '
$c = $ws.Range('D9')
$c.Value = "'0.536"
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('E10').Value = '  +0.90%  '
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('E12').Value = '  -0.71%  '
$c = $ws.Range('D13')
$c.Value = "'37.25"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -3.01%  '
$c = $ws.Range('D14')
$c.Value = "'0.0000245"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('D15').Value = '4.431.60'
$ws.Range('E15').Value = '  +1.20%  '
$ws.Range('D16').Value = '3.785.31'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('D17').Value = '69.238.54'
$ws.Range('E17').Value = '  +0.24%  '
$c = $ws.Range('D18')
$c.Value = "'7.40"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('E19').Value = '  -0.31%  '
$c = $ws.Range('D20')
$c.Value = "'17.29"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +1.04%  '
$c = $ws.Range('D21')
$c.Value = "'11.33"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +4.37%  '
$c = $ws.Range('D22')
$c.Value = "'488.41"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.16%  '
$c = $ws.Range('D23')
$c.Value = "'0.722"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('E24').Value = '  -2.17%  '
$c = $ws.Range('D25')
$c.Value = "'84.58"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.03%  '
$c = $ws.Range('D26')
$c.Value = "'2.25"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -3.60%  '
$c = $ws.Range('D27')
$c.Value = "'12.20"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.23%  '
$c = $ws.Range('D28')
$c.Value = "'10.05"
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -2.49%  '
$ws.Range('E29').Value = '  +0.01%  '
$c = $ws.Range('D30')
$c.Value = "'2.96"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.79%  '
$c = $ws.Range('D31')
$c.Value = "'8.04"
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.24%  '
$c = $ws.Range('D32')
$c.Value = "'2.39"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -5.72%  '
$ws.Range('D33').Value = '3.947.88'
$ws.Range('E33').Value = '  +1.26%  '
$c = $ws.Range('D34')
$c.Value = "'31.82"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').Value = '3.743.63'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  +5.57%  '
$ws.Range('E38').Value = '  +0.16%  '
$c = $ws.Range('D39')
$c.Value = "'5.91"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('E40').Value = '  -0.04%  '
$c = $ws.Range('D41')
$c.Value = "'0.322"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.44%  '
$c = $ws.Range('D42')
$c.Value = "'3.02"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.42%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D43')
$c.Value = "'48.58"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D44')
$c.Value = "'1.99"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.55%  '
$c = $ws.Range('D45')
$c.Value = "'421.23"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -3.91%  '
$c = $ws.Range('D47')
$c.Value = "'8.37"
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D48').Value = '2.828.03'
$ws.Range('E48').Value = '  +0.94%  '
$c = $ws.Range('D49')
$c.Value = "'39.69"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.68%  '
$c = $ws.Range('D50')
$c.Value = "'140.82"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -0.17%  '
$ws.Range('E51').Value = '  +4.39%  '
